$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1785714285714286
$ws.Range("C2").Value = 0.6038961038961039
$ws.Range("J2").Value = 0.01623376623376623
$ws.Range("P2").Value = 0.1233766233766234
$ws.Range("S2").Value = 0.07792207792207792
$ws.Range("B3").Value = 0.01595744680851064
$ws.Range("C3").Value = 0.02659574468085106
$ws.Range("J3").Value = 0.03723404255319149
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.1702127659574468
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.06273062730627306
$ws.Range("D6").Value = 0.01107011070110701
$ws.Range("F6").Value = 0.1254612546125461
$ws.Range("J6").Value = 0.2066420664206642
$ws.Range("O6").Value = 0.02952029520295203
$ws.Range("Q6").Value = 0.1955719557195572
$ws.Range("R6").Value = 0.05535055350553506
$ws.Range("S6").Value = 0.3136531365313653
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.0101010101010101
$ws.Range("F7").Value = 0.07575757575757576
$ws.Range("J7").Value = 0.1212121212121212
$ws.Range("O7").Value = 0.03535353535353535
$ws.Range("Q7").Value = 0.2070707070707071
$ws.Range("R7").Value = 0.1060606060606061
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.08454106280193237
$ws.Range("D8").Value = 0.01690821256038647
$ws.Range("E8").Value = 0.002415458937198068
$ws.Range("F8").Value = 0.07729468599033816
$ws.Range("J8").Value = 0.1038647342995169
$ws.Range("O8").Value = 0.04106280193236715
$ws.Range("Q8").Value = 0.2053140096618357
$ws.Range("R8").Value = 0.06763285024154589
$ws.Range("S8").Value = 0.4009661835748792
$ws.Range("B9").Value = 0.09625668449197861
$ws.Range("D9").Value = 0.0053475935828877
$ws.Range("F9").Value = 0.1229946524064171
$ws.Range("J9").Value = 0.08021390374331551
$ws.Range("O9").Value = 0.0427807486631016
$ws.Range("Q9").Value = 0.160427807486631
$ws.Range("R9").Value = 0.09090909090909091
$ws.Range("S9").Value = 0.4010695187165775
$ws.Range("B10").Value = 0.1342905405405405
$ws.Range("D10").Value = 0.02280405405405405
$ws.Range("E10").Value = 0.0008445945945945946
$ws.Range("F10").Value = 0.08108108108108109
$ws.Range("J10").Value = 0.09797297297297297
$ws.Range("O10").Value = 0.02702702702702703
$ws.Range("Q10").Value = 0.1925675675675676
$ws.Range("R10").Value = 0.07094594594594594
$ws.Range("S10").Value = 0.3724662162162162
$ws.Range("G11").Value = 0.1184210526315789
$ws.Range("J11").Value = 0.1118421052631579
$ws.Range("K11").Value = 0.194078947368421
$ws.Range("L11").Value = 0.5657894736842105
$ws.Range("S11").Value = 0.009868421052631578
$ws.Range("G12").Value = 0.7175141242937854
$ws.Range("J12").Value = 0.2146892655367232
$ws.Range("K12").Value = 0.005649717514124294
$ws.Range("L12").Value = 0.02259887005649718
$ws.Range("S12").Value = 0.03954802259887006
$ws.Range("G13").Value = 0.6833333333333333
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.01666666666666667
$ws.Range("F15").Value = 0.02016129032258064
$ws.Range("H15").Value = 0.1370967741935484
$ws.Range("I15").Value = 0.07661290322580645
$ws.Range("J15").Value = 0.3225806451612903
$ws.Range("K15").Value = 0.07661290322580645
$ws.Range("O15").Value = 0.08064516129032258
$ws.Range("S15").Value = 0.2862903225806452
$ws.Range("F16").Value = 0.03045685279187817
$ws.Range("H16").Value = 0.1979695431472081
$ws.Range("I16").Value = 0.1015228426395939
$ws.Range("J16").Value = 0.3807106598984771
$ws.Range("K16").Value = 0.06091370558375635
$ws.Range("M16").Value = 0.02538071065989848
$ws.Range("N16").Value = 0.005076142131979695
$ws.Range("O16").Value = 0.06091370558375635
$ws.Range("S16").Value = 0.1370558375634518
$ws.Range("F17").Value = 0.02293577981651376
$ws.Range("H17").Value = 0.1697247706422018
$ws.Range("I17").Value = 0.07339449541284404
$ws.Range("J17").Value = 0.4105504587155963
$ws.Range("K17").Value = 0.1123853211009174
$ws.Range("M17").Value = 0.01605504587155963
$ws.Range("O17").Value = 0.05963302752293578
$ws.Range("S17").Value = 0.1353211009174312
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.1566265060240964
$ws.Range("I18").Value = 0.0783132530120482
$ws.Range("J18").Value = 0.4518072289156627
$ws.Range("K18").Value = 0.08433734939759036
$ws.Range("M18").Value = 0.03012048192771084
$ws.Range("N18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.03614457831325301
$ws.Range("S18").Value = 0.144578313253012
$ws.Range("F19").Value = 0.02369281045751634
$ws.Range("H19").Value = 0.1977124183006536
$ws.Range("I19").Value = 0.08496732026143791
$ws.Range("J19").Value = 0.3513071895424836
$ws.Range("K19").Value = 0.1209150326797386
$ws.Range("M19").Value = 0.03594771241830065
$ws.Range("N19").Value = 0.001633986928104575
$ws.Range("O19").Value = 0.06535947712418301
$ws.Range("S19").Value = 0.1184640522875817
